$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This report tracks localization hand-off status for a single source file.
# The file being tracked changed from UUID 8ad3246c-... to b324597a-...,
# a new hand-off round was generated, and the (not-yet-existing) hand-back
# info for the new file was cleared out.
# ---------------------------------------------------------------------------

$oldGuid = "8ad3246c-fe12-4e35-bf27-126707fad4b0"
$newGuid = "b324597a-54a8-4c83-899f-154686303bcd"

# =============================== Overview ===================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newGuid.md"
$ov.Range("B2").Value = "e2e\$newGuid.md"
$ov.Range("G2").Value = "2016-08-25 13:02:40"

# The hyperlink target (rId2 -> the GitHub blob URL) does not change, only
# the display text that mirrors the new file name. Rebuild the single
# hyperlink on this sheet so the display text is refreshed in place.
$ovAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a10a33a06c0782d2d861f3ec222c9bb19ea029c/e2e/$oldGuid.md"
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), $ovAddr, "", "", "e2e\$newGuid.md") | Out-Null

# =============================== zh-cn =======================================
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$newGuid.md"
$zh.Range("G2").Value = "$newGuid.e902cda3b0f0dc8e97008c492cad48dc21cb177a.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-25 13:02:35"

# Latest Target File / Latest Handback File / Latest Handback DateTime reset
# for the new file (no hand-back has happened for it yet).
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

# Only the A2 hyperlink remains; the I2 hyperlink (to the hand-back target)
# is removed entirely since there is no hand-back target for the new file.
$zhAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a10a33a06c0782d2d861f3ec222c9bb19ea029c/e2e/$oldGuid.md"
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhAddr, "", "", "$newGuid.md") | Out-Null

$zh.Columns.Item(9).ColumnWidth = 17.8
$zh.Columns.Item(10).ColumnWidth = 20.8

# =============================== de-de =======================================
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$newGuid.md"
$de.Range("G2").Value = "$newGuid.e902cda3b0f0dc8e97008c492cad48dc21cb177a.de-de.xlf"
$de.Range("H2").Value = "2016-08-25 13:02:40"

$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$deAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a10a33a06c0782d2d861f3ec222c9bb19ea029c/e2e/$oldGuid.md"
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deAddr, "", "", "$newGuid.md") | Out-Null

$de.Columns.Item(9).ColumnWidth = 17.8
$de.Columns.Item(10).ColumnWidth = 20.8

Write-Host "Localization status report refreshed for handoff."
